$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 hold the per-game "Date" column, currently mis-formatted as
# "4-26-2013-14" (day-month-season mash-up). Correct it to the real game
# date in ISO form, "2014-04-26". A leading apostrophe forces the engine
# to keep it as literal text instead of re-parsing the ISO-like string
# back into a serial date value.
for ($i = 2; $i -le 31; $i++) {
    $ws.Cells.Item($i, 58).Value = "'2014-04-26"
}
